# Apply the edits described by the commit diff:
#  1. Resize the saved window dimensions recorded for the workbook.
#  2. Change the active selection on Sheet1 from C11 to A2 (with A2:A7 selected).
#  3. Change the numeric values in A2:A7 from 2 to 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Window size (best effort - mirrors the bookViews/workbookView change in workbook.xml)
try {
    $excel.ActiveWindow.Width = 30240
    $excel.ActiveWindow.Height = 11500
} catch {
    # Not fatal if the hosting environment does not allow resizing the window.
}

# 2) Update the data values in column A (rows 2-7) from 2 to 7
$ws.Range("A2:A7").Value = 7

# 3) Update the selection to A2:A7, with A2 as the active cell
$ws.Range("A2:A7").Select()
